# Updates cryptos list price/volume data (and the Aave/Bittensor row swap)
# per the scraped coinranking.com snapshot for this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    # Force text storage so numeric-looking strings (e.g. "547.60",
    # "135.30") keep trailing zeros / dotted-thousands formatting
    # instead of being auto-coerced to a Double by COM.
    $c.NumberFormat = "@"
    $c.Value = $text
    # Restore the default "Normal" style so we do not leave the cell
    # pinned to the Text number format (matches original styling).
    $c.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '59.013.08'
Set-TextCell 2 5 '  -2.40%  '

# Row 3
Set-TextCell 3 4 '2.546.81'
Set-TextCell 3 5 '  -2.88%  '

# Row 4
Set-TextCell 4 5 '  +0.00%  '

# Row 5
Set-TextCell 5 4 '547.60'
Set-TextCell 5 5 '  -4.18%  '

# Row 6
Set-TextCell 6 4 '139.55'
Set-TextCell 6 5 '  -4.32%  '

# Row 7
Set-TextCell 7 5 '  +0.35%  '

# Row 8
Set-TextCell 8 4 '0.588'
Set-TextCell 8 5 '  -2.09%  '

# Row 9
Set-TextCell 9 4 '2.546.15'
Set-TextCell 9 5 '  -3.59%  '

# Row 10
Set-TextCell 10 4 '6.66'
Set-TextCell 10 5 '  -1.35%  '

# Row 11
Set-TextCell 11 5 '  -2.14%  '

# Row 12
Set-TextCell 12 5 '  +5.62%  '

# Row 13
Set-TextCell 13 4 '0.351'
Set-TextCell 13 5 '  +1.83%  '

# Row 14
Set-TextCell 14 4 '2.995.60'
Set-TextCell 14 5 '  -2.83%  '

# Row 15
Set-TextCell 15 4 '59.039.74'
Set-TextCell 15 5 '  -2.27%  '

# Row 16
Set-TextCell 16 4 '22.97'
Set-TextCell 16 5 '  +3.87%  '

# Row 17
Set-TextCell 17 5 '  -1.90%  '

# Row 18
Set-TextCell 18 4 '2.561.95'
Set-TextCell 18 5 '  -2.61%  '

# Row 19
Set-TextCell 19 5 '  -0.68%  '

# Row 20
Set-TextCell 20 4 '335.03'
Set-TextCell 20 5 '  -2.03%  '

# Row 21
Set-TextCell 21 4 '10.19'
Set-TextCell 21 5 '  -2.17%  '

# Row 22
Set-TextCell 22 4 '6.36'
Set-TextCell 22 5 '  -0.03%  '

# Row 23
Set-TextCell 23 4 '0.991'
Set-TextCell 23 5 '  -0.83%  '

# Row 24
Set-TextCell 24 5 '  +5.14%  '

# Row 25
Set-TextCell 25 4 '62.38'
Set-TextCell 25 5 '  -5.26%  '

# Row 26
Set-TextCell 26 5 '  +0.67%  '

# Row 27
Set-TextCell 27 5 '  -3.72%  '

# Row 28
Set-TextCell 28 4 '7.33'
Set-TextCell 28 5 '  -0.57%  '

# Row 29
Set-TextCell 29 4 '0.0₃0754'
Set-TextCell 29 5 '  -5.39%  '

# Row 30
Set-TextCell 30 5 '  +0.16%  '

# Row 31
Set-TextCell 31 5 '  -0.58%  '

# Row 32
Set-TextCell 32 5 '  -3.23%  '

# Row 33
Set-TextCell 33 4 '158.10'
Set-TextCell 33 5 '  -1.13%  '

# Row 34
Set-TextCell 34 4 '18.88'
Set-TextCell 34 5 '  -1.53%  '

# Row 35
Set-TextCell 35 4 '4.03'
Set-TextCell 35 5 '  -1.87%  '

# Row 36
Set-TextCell 36 4 '1.15'
Set-TextCell 36 5 '  +0.11%  '

# Row 37
Set-TextCell 37 5 '  -1.22%  '

# Row 38
Set-TextCell 38 4 '37.27'
Set-TextCell 38 5 '  -0.76%  '

# Row 39
Set-TextCell 39 4 '0.837'
Set-TextCell 39 5 '  -5.51%  '

# Row 40
Set-TextCell 40 5 '  -4.00%  '

# Row 41
Set-TextCell 41 4 '3.63'
Set-TextCell 41 5 '  -0.87%  '

# Row 42
Set-TextCell 42 2 'Aave'
Set-TextCell 42 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 42 4 '135.30'
Set-TextCell 42 5 '  +6.71%  '

# Row 43
Set-TextCell 43 2 'Bittensor'
Set-TextCell 43 3 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 43 4 '282.04'
Set-TextCell 43 5 '  -5.33%  '

# Row 44
Set-TextCell 44 4 '0.999'
Set-TextCell 44 5 '  +0.39%  '

# Row 45
Set-TextCell 45 5 '  -1.80%  '

# Row 46
Set-TextCell 46 4 '10.66'
Set-TextCell 46 5 '  -0.07%  '

# Row 47
Set-TextCell 47 5 '  -3.23%  '

# Row 48
Set-TextCell 48 5 '  -3.11%  '

# Row 49
Set-TextCell 49 4 '0.0231'
Set-TextCell 49 5 '  -2.29%  '

# Row 50
Set-TextCell 50 4 '1.947.22'
Set-TextCell 50 5 '  -0.77%  '

# Row 51
Set-TextCell 51 4 '18.36'
Set-TextCell 51 5 '  -2.03%  '

